# Apply cryptocurrency price/volume updates for Thu Jul 13 05:35:05 UTC 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force columns D and E to Text format so that numeric-looking
# strings (e.g. "0.9999") are not coerced into numbers by Excel's type inference.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.305.59"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.866.17"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "243.45"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4724"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").Value = "0.2870"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").Value = "0.06472"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "0.07791"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "97.02"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.867.71"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "0.7199"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "5.141"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "280.29"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "30.299.64"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "13.00"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("D20").Value = "0.000007459"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "2.108.74"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "5.227"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "6.248"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "162.10"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "8.987"
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").Value = "18.68"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "1.876"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").Value = "0.09629"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "4.215"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("D33").Value = "4.110"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").Value = "0.04773"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "0.6832"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "0.01888"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "2.835"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "75.13"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "6.197"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "1.936"
$ws.Range("E42").Value = "  -4.84%  "
$ws.Range("D43").Value = "0.4204"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "0.9987"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "0.8267"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "100.47"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "9.596"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "6.950"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("D50").Value = "0.05761"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "885.72"
$ws.Range("E51").Value = "  -3.41%  "

# Restore the original (default) cell style now that the text values are set,
# so the cells keep matching the workbook's original (unstyled) appearance.
$dataRange.Style = "Normal"

